$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028396595994657
$ws.Range("D2").Value = 1.037503117341549
$ws.Range("E2").Value = 1.028371722273546
$ws.Range("F2").Value = 1.046158541185202
$ws.Range("I2").Value = 1.036782525493392
$ws.Range("J2").Value = 1.03354914486589
$ws.Range("K2").Value = 1.040293816227196
$ws.Range("L2").Value = 1.031188732285669
$ws.Range("M2").Value = 1.048924760530941
$ws.Range("N2").Value = 1.015204788863329
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029173548958836
$ws.Range("D3").Value = 1.038127352388706
$ws.Range("E3").Value = 1.029026977326274
$ws.Range("F3").Value = 1.047015410621764
$ws.Range("I3").Value = 1.036969898359139
$ws.Range("J3").Value = 1.03396757382998
$ws.Range("K3").Value = 1.040728363716059
$ws.Range("L3").Value = 1.031652320913652
$ws.Range("M3").Value = 1.049593098310582
$ws.Range("N3").Value = 1.015343802872476
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02967686913955
$ws.Range("D4").Value = 1.038531741333731
$ws.Range("E4").Value = 1.029451863654791
$ws.Range("F4").Value = 1.047570886002336
$ws.Range("I4").Value = 1.037090052961926
$ws.Range("J4").Value = 1.034238226559727
$ws.Range("K4").Value = 1.041009320332642
$ws.Range("L4").Value = 1.031952501493061
$ws.Range("M4").Value = 1.050025937035922
$ws.Range("N4").Value = 1.01543370115611
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029888601902568
$ws.Range("D5").Value = 1.038701856400068
$ws.Range("E5").Value = 1.029630697621868
$ws.Range("F5").Value = 1.047804650822118
$ws.Range("I5").Value = 1.03714030459179
$ws.Range("J5").Value = 1.03435198403223
$ws.Range("K5").Value = 1.041127379325015
$ws.Range("L5").Value = 1.032078745728962
$ws.Range("M5").Value = 1.050207991503741
$ws.Range("N5").Value = 1.015471481231301
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029924160743819
$ws.Range("D6").Value = 1.038730425861348
$ws.Range("E6").Value = 1.02966073699363
$ws.Range("F6").Value = 1.047843915151596
$ws.Range("I6").Value = 1.037148726710147
$ws.Range("J6").Value = 1.034371082912848
$ws.Range("K6").Value = 1.04114719867392
$ws.Range("L6").Value = 1.032099945484094
$ws.Range("M6").Value = 1.050238564421085
$ws.Range("N6").Value = 1.015477823886975
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029679697787208
$ws.Range("D7").Value = 1.038534013989894
$ws.Range("E7").Value = 1.029454252413369
$ws.Range("F7").Value = 1.047574008627196
$ws.Range("I7").Value = 1.037090725454718
$ws.Range("J7").Value = 1.034239746692234
$ws.Range("K7").Value = 1.041010898061192
$ws.Range("L7").Value = 1.031954188186631
$ws.Range("M7").Value = 1.050028369308895
$ws.Range("N7").Value = 1.015434206027621
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028659049773856
$ws.Range("D8").Value = 1.037713982648078
$ws.Range("E8").Value = 1.02859298284202
$ws.Range("F8").Value = 1.046447911042402
$ws.Range("I8").Value = 1.036846073829113
$ws.Range("J8").Value = 1.033690574576061
$ws.Range("K8").Value = 1.04044071927219
$ws.Range("L8").Value = 1.031345360370556
$ws.Range("M8").Value = 1.049150548763816
$ws.Range("N8").Value = 1.015251780013753
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026865057268223
$ws.Range("D9").Value = 1.036272643238473
$ws.Range("E9").Value = 1.027082235645842
$ws.Range("F9").Value = 1.04447151408279
$ws.Range("I9").Value = 1.036406670139277
$ws.Range("J9").Value = 1.032722163307501
$ws.Range("K9").Value = 1.039434331875316
$ws.Range("L9").Value = 1.030274181407679
$ws.Range("M9").Value = 1.047606702103299
$ws.Range("N9").Value = 1.014929935750306
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025672209217339
$ws.Range("D10").Value = 1.035314323523389
$ws.Range("E10").Value = 1.026079837437543
$ws.Range("F10").Value = 1.043159372078729
$ws.Range("I10").Value = 1.036108207432216
$ws.Range("J10").Value = 1.032076161905519
$ws.Range("K10").Value = 1.03876237372915
$ws.Range("L10").Value = 1.029561258083202
$ws.Range("M10").Value = 1.046579584576933
$ws.Range("N10").Value = 1.014715140411139
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025156462081797
$ws.Range("D11").Value = 1.034899994555398
$ws.Range("E11").Value = 1.025646942692955
$ws.Range("F11").Value = 1.042592520338202
$ws.Range("I11").Value = 1.035977669632167
$ws.Range("J11").Value = 1.031796358933096
$ws.Range("K11").Value = 1.03847118080998
$ws.Range("L11").Value = 1.02925285592062
$ws.Range("M11").Value = 1.04613535515377
$ws.Range("N11").Value = 1.014622082269961
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024965007247248
$ws.Range("D12").Value = 1.034746190872969
$ws.Range("E12").Value = 1.025486321099893
$ws.Range("F12").Value = 1.042382165958545
$ws.Range("I12").Value = 1.035928987400708
$ws.Range("J12").Value = 1.03169241692192
$ws.Range("K12").Value = 1.038362985553196
$ws.Range("L12").Value = 1.029138347893506
$ws.Range("M12").Value = 1.045970428395712
$ws.Range("N12").Value = 1.014587509234082
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025006069652503
$ws.Range("D13").Value = 1.034779177885244
$ws.Range("E13").Value = 1.025520767035641
$ws.Range("F13").Value = 1.042427278620861
$ws.Range("I13").Value = 1.035939438700616
$ws.Range("J13").Value = 1.031714713300717
$ws.Range("K13").Value = 1.038386195284744
$ws.Range("L13").Value = 1.029162908137137
$ws.Range("M13").Value = 1.04600580210355
$ws.Range("N13").Value = 1.014594925581425
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025140633975878
$ws.Range("D14").Value = 1.034887279114435
$ws.Range("E14").Value = 1.025633662084926
$ws.Range("F14").Value = 1.042575128295599
$ws.Range("I14").Value = 1.035973649510216
$ws.Range("J14").Value = 1.03178776726535
$ws.Range("K14").Value = 1.038462238026292
$ws.Range("L14").Value = 1.029243389699182
$ws.Range("M14").Value = 1.046121720617498
$ws.Range("N14").Value = 1.01461922459416
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02522355900445
$ws.Range("D15").Value = 1.034953896700107
$ws.Range("E15").Value = 1.025703243671014
$ws.Range("F15").Value = 1.042666249871686
$ws.Range("I15").Value = 1.035994702147582
$ws.Range("J15").Value = 1.031832776841381
$ws.Range("K15").Value = 1.038509086106307
$ws.Range("L15").Value = 1.029292983219892
$ws.Range("M15").Value = 1.046193152472238
$ws.Range("N15").Value = 1.014634195091092
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025706453908702
$ws.Range("D16").Value = 1.035341834605841
$ws.Range("E16").Value = 1.026108591639754
$ws.Range("F16").Value = 1.043197020035969
$ws.Range("I16").Value = 1.036116843448148
$ws.Range("J16").Value = 1.032094729931085
$ws.Range("K16").Value = 1.038781694505488
$ws.Range("L16").Value = 1.02958173212816
$ws.Range("M16").Value = 1.046609077715825
$ws.Range("N16").Value = 1.014721315352506
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026009566952847
$ws.Range("D17").Value = 1.035585347778309
$ws.Range("E17").Value = 1.026363165070512
$ws.Range("F17").Value = 1.043530311538672
$ws.Range("I17").Value = 1.036193111617565
$ws.Range("J17").Value = 1.032259025659779
$ws.Range("K17").Value = 1.038952633804007
$ws.Range("L17").Value = 1.029762937585062
$ws.Range("M17").Value = 1.046870116932763
$ws.Range("N17").Value = 1.014775950440941
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026186441249003
$ws.Range("D18").Value = 1.035727445528833
$ws.Range("E18").Value = 1.026511764300405
$ws.Range("F18").Value = 1.043724841512254
$ws.Range("I18").Value = 1.036237471919451
$ws.Range("J18").Value = 1.032354848762454
$ws.Range("K18").Value = 1.039052317460319
$ws.Range("L18").Value = 1.029868660351143
$ws.Range("M18").Value = 1.047022426669137
$ws.Range("N18").Value = 1.014807813275299
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026246763238564
$ws.Range("D19").Value = 1.035775907416085
$ws.Range("E19").Value = 1.02656245154052
$ws.Range("F19").Value = 1.043791192628538
$ws.Range("I19").Value = 1.036252576300917
$ws.Range("J19").Value = 1.032387520598677
$ws.Range("K19").Value = 1.0390863031783
$ws.Range("L19").Value = 1.029904713904377
$ws.Range("M19").Value = 1.047074368774932
$ws.Range("N19").Value = 1.014818676835113
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025977038167238
$ws.Range("D20").Value = 1.035559214833639
$ws.Range("E20").Value = 1.026335840270205
$ws.Range("F20").Value = 1.043494539396378
$ws.Range("I20").Value = 1.036184941753238
$ws.Range("J20").Value = 1.032241399082019
$ws.Range("K20").Value = 1.038934295936086
$ws.Range("L20").Value = 1.029743492976603
$ws.Range("M20").Value = 1.046842104712705
$ws.Range("N20").Value = 1.01477008911529
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025101004902626
$ws.Range("D21").Value = 1.034855443316254
$ws.Range("E21").Value = 1.025600412465887
$ws.Range("F21").Value = 1.042531584724592
$ws.Range("I21").Value = 1.035963580646754
$ws.Range("J21").Value = 1.031766254979282
$ws.Range("K21").Value = 1.038439846240566
$ws.Range("L21").Value = 1.0292196885963
$ws.Range("M21").Value = 1.046087583286061
$ws.Range("N21").Value = 1.014612069332783
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02455088363179
$ws.Range("D22").Value = 1.034413513756514
$ws.Range("E22").Value = 1.025139032068346
$ws.Range("F22").Value = 1.041927292752514
$ws.Range("I22").Value = 1.035823276195466
$ws.Range("J22").Value = 1.031467451594285
$ws.Range("K22").Value = 1.03812877446343
$ws.Range("L22").Value = 1.02889062039701
$ws.Range("M22").Value = 1.04561364848406
$ws.Range("N22").Value = 1.014512675181517
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024842448670596
$ws.Range("D23").Value = 1.034647735299219
$ws.Range("E23").Value = 1.025383521939812
$ws.Range("F23").Value = 1.042247529091144
$ws.Range("I23").Value = 1.035897760680864
$ws.Range("J23").Value = 1.031625858309354
$ws.Range("K23").Value = 1.038293697193525
$ws.Range("L23").Value = 1.029065039751972
$ws.Range("M23").Value = 1.045864845817828
$ws.Range("N23").Value = 1.014565369618154
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025991736298136
$ws.Range("D24").Value = 1.035571022999352
$ws.Range("E24").Value = 1.026348186829023
$ws.Range("F24").Value = 1.043510702896543
$ws.Range("I24").Value = 1.036188633751544
$ws.Range("J24").Value = 1.032249363798101
$ws.Range("K24").Value = 1.038942582098476
$ws.Range("L24").Value = 1.029752279069728
$ws.Range("M24").Value = 1.046854762074988
$ws.Range("N24").Value = 1.014772737611218
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.0273283005845
$ws.Range("D25").Value = 1.036644818668566
$ws.Range("E25").Value = 1.027471968203913
$ws.Range("F25").Value = 1.044981507474755
$ws.Range("I25").Value = 1.036521244421398
$ws.Range("J25").Value = 1.032972595866435
$ws.Range("K25").Value = 1.039694695292042
$ws.Range("L25").Value = 1.030550902210649
$ws.Range("M25").Value = 1.048005458728663
$ws.Range("N25").Value = 1.015013183178443
